# "add some test sample"
#
# The sheet has three copies of a small report table (rows 1-8, 11-18,
# 21-28, 30-37), each with a label column (A) and a "测试一" value
# column (B). This change adds a third column ("测试二") of results
# next to the second table (rows 11-18) only.
#
# New unique shared strings are introduced by this edit, in this exact
# order, to match the target sharedStrings.xml append order:
#   17 "40W样本中10W样本"  (C12 - sub-title, mirrors B12)
#   18 "48684.8 s"          (C18 - timing result, mirrors B18)
#   19 "测试二"             (C11 - column header, mirrors B11)
# so the new values are written in that same order below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C12").Value = "40W样本中10W样本"
$ws.Range("C18").Value = "48684.8 s"
$ws.Range("C11").Value = "测试二"
$ws.Range("C13").Value = 139181
$ws.Range("C14").Value = 100000
$ws.Range("C15").Value = 220487
$ws.Range("C16").Value = 7
$ws.Range("C17").Value = 5

# Give the new column the same look as column B (header style, wrap
# style on the sub-title row, right-aligned numbers, etc.)
$ws.Range("B11:B18").Copy()
$ws.Range("C11:C18").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Scroll/select roughly where the new data was entered.
$ws.Range("F20").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
